$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 12,20
$data[0,0] = "ECs"
$data[0,1] = "Wnt2"
$data[0,2] = "Fzd9"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.01070233333333333
$data[0,7] = 0.032107
$data[0,8] = 0.006017198313602724
$data[0,9] = 0.006017198313602724
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.514018
$data[0,13] = 1.542054
$data[0,14] = 0.2743122022496015
$data[0,15] = 0.2743122022496015
$data[0,16] = 0.005501191975333333
$data[0,17] = 0.049510727778
$data[0,18] = 0.001650590920776952
$data[0,19] = 0.001650590920776951

$data[1,0] = "ECs"
$data[1,1] = "Wnt2"
$data[1,2] = "Fzd9"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.01070233333333333
$data[1,7] = 0.032107
$data[1,8] = 0.006017198313602724
$data[1,9] = 0.006017198313602724
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.5446803333333333
$data[1,13] = 1.634041
$data[1,14] = 0.2906755439667749
$data[1,15] = 0.2906755439667749
$data[1,16] = 0.005829350487444444
$data[1,17] = 0.052464154387
$data[1,18] = 0.001749052392962432
$data[1,19] = 0.001749052392962432

$data[2,0] = "ECs"
$data[2,1] = "Wnt2"
$data[2,2] = "Fzd9"
$data[2,3] = "M1"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.01070233333333333
$data[2,7] = 0.032107
$data[2,8] = 0.006017198313602724
$data[2,9] = 0.006017198313602724
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.1138553333333334
$data[2,13] = 0.341566
$data[2,14] = 0.0607603376234473
$data[2,15] = 0.06076033762344728
$data[2,16] = 0.001218517729111111
$data[2,17] = 0.010966659562
$data[2,18] = 0.0003656070010817392
$data[2,19] = 0.0003656070010817391

$data[3,0] = "ECs"
$data[3,1] = "Wnt2"
$data[3,2] = "Fzd9"
$data[3,3] = "M2"
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.01070233333333333
$data[3,7] = 0.032107
$data[3,8] = 0.006017198313602724
$data[3,9] = 0.006017198313602724
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.1469773333333333
$data[3,13] = 0.440932
$data[3,14] = 0.07843631154442145
$data[3,15] = 0.07843631154442145
$data[3,16] = 0.001573000413777778
$data[3,17] = 0.014157003724
$data[3,18] = 0.0004719668415503106
$data[3,19] = 0.0004719668415503106

$data[4,0] = "ECs"
$data[4,1] = "Wnt2"
$data[4,2] = "Fzd9"
$data[4,3] = "Neutro"
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.01070233333333333
$data[4,7] = 0.032107
$data[4,8] = 0.006017198313602724
$data[4,9] = 0.006017198313602724
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.04771466666666666
$data[4,13] = 0.143144
$data[4,14] = 0.02546353492083737
$data[4,15] = 0.02546353492083737
$data[4,16] = 0.0005106582675555555
$data[4,17] = 0.004595924408
$data[4,18] = 0.0001532191393840267
$data[4,19] = 0.0001532191393840267

$data[5,0] = "ECs"
$data[5,1] = "Wnt2"
$data[5,2] = "Fzd9"
$data[5,3] = "sCs"
$data[5,4] = 1
$data[5,5] = 0.3333333333333333
$data[5,6] = 0.01070233333333333
$data[5,7] = 0.032107
$data[5,8] = 0.006017198313602724
$data[5,9] = 0.006017198313602724
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.5065973333333333
$data[5,13] = 1.519792
$data[5,14] = 0.2703520696949175
$data[5,15] = 0.2703520696949175
$data[5,16] = 0.005421773527111111
$data[5,17] = 0.048795961744
$data[5,18] = 0.001626762017847264
$data[5,19] = 0.001626762017847264

$data[6,0] = "FAPs"
$data[6,1] = "Wnt2"
$data[6,2] = "Fzd9"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1.767921666666667
$data[6,7] = 5.303765
$data[6,8] = 0.9939828016863973
$data[6,9] = 0.9939828016863973
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.514018
$data[6,13] = 1.542054
$data[6,14] = 0.2743122022496015
$data[6,15] = 0.2743122022496015
$data[6,16] = 0.9087435592566667
$data[6,17] = 8.17869203331
$data[6,18] = 0.2726616113288245
$data[6,19] = 0.2726616113288245

$data[7,0] = "FAPs"
$data[7,1] = "Wnt2"
$data[7,2] = "Fzd9"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.767921666666667
$data[7,7] = 5.303765
$data[7,8] = 0.9939828016863973
$data[7,9] = 0.9939828016863973
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.5446803333333333
$data[7,13] = 1.634041
$data[7,14] = 0.2906755439667749
$data[7,15] = 0.2906755439667749
$data[7,16] = 0.9629521627072223
$data[7,17] = 8.666569464365
$data[7,18] = 0.2889264915738124
$data[7,19] = 0.2889264915738124

$data[8,0] = "FAPs"
$data[8,1] = "Wnt2"
$data[8,2] = "Fzd9"
$data[8,3] = "M1"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 1.767921666666667
$data[8,7] = 5.303765
$data[8,8] = 0.9939828016863973
$data[8,9] = 0.9939828016863973
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.1138553333333334
$data[8,13] = 0.341566
$data[8,14] = 0.0607603376234473
$data[8,15] = 0.06076033762344728
$data[8,16] = 0.2012873106655556
$data[8,17] = 1.81158579599
$data[8,18] = 0.06039473062236556
$data[8,19] = 0.06039473062236555

$data[9,0] = "FAPs"
$data[9,1] = "Wnt2"
$data[9,2] = "Fzd9"
$data[9,3] = "M2"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 1.767921666666667
$data[9,7] = 5.303765
$data[9,8] = 0.9939828016863973
$data[9,9] = 0.9939828016863973
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 0.1469773333333333
$data[9,13] = 0.440932
$data[9,14] = 0.07843631154442145
$data[9,15] = 0.07843631154442145
$data[9,16] = 0.2598444121088889
$data[9,17] = 2.33859970898
$data[9,18] = 0.07796434470287114
$data[9,19] = 0.07796434470287114

$data[10,0] = "FAPs"
$data[10,1] = "Wnt2"
$data[10,2] = "Fzd9"
$data[10,3] = "Neutro"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 1.767921666666667
$data[10,7] = 5.303765
$data[10,8] = 0.9939828016863973
$data[10,9] = 0.9939828016863973
$data[10,10] = 2
$data[10,11] = 0.6666666666666666
$data[10,12] = 0.04771466666666666
$data[10,13] = 0.143144
$data[10,14] = 0.02546353492083737
$data[10,15] = 0.02546353492083737
$data[10,16] = 0.08435579301777778
$data[10,17] = 0.75920213716
$data[10,18] = 0.02531031578145335
$data[10,19] = 0.02531031578145335

$data[11,0] = "FAPs"
$data[11,1] = "Wnt2"
$data[11,2] = "Fzd9"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 1.767921666666667
$data[11,7] = 5.303765
$data[11,8] = 0.9939828016863973
$data[11,9] = 0.9939828016863973
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.5065973333333333
$data[11,13] = 1.519792
$data[11,14] = 0.2703520696949175
$data[11,15] = 0.2703520696949175
$data[11,16] = 0.8956244018755557
$data[11,17] = 8.06061961688
$data[11,18] = 0.2687253076770703
$data[11,19] = 0.2687253076770702

$ws.Range("A2:T13").Value = $data